$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 19) to the stock/news-tone table, continuing
# the existing daily rows (row 18 was the last one, for 2025-07-18).

$row = 19

$ws.Cells.Item($row, 1).Value = 45859
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = 6305.60009765625
$ws.Cells.Item($row, 3).Value = 6336.080078125
$ws.Cells.Item($row, 4).Value = 6303.7900390625
$ws.Cells.Item($row, 5).Value = 6304.740234375
$ws.Cells.Item($row, 6).Value = 5010840000
$ws.Cells.Item($row, 7).Value = 0.0013991348828683
$ws.Cells.Item($row, 8).Value = 1

# daily_headlines (I) has no news for this date -> blank cell.
# daily_headlines_clean (J) mirrors the Python "nan" placeholder text.
$ws.Cells.Item($row, 10).Value = "nan"

$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 1
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = 0
$ws.Cells.Item($row, 19).Value = 0
$ws.Cells.Item($row, 20).Value = 0
$ws.Cells.Item($row, 21).Value = 0
$ws.Cells.Item($row, 22).Value = 0
$ws.Cells.Item($row, 23).Value = 0
$ws.Cells.Item($row, 24).Value = 0
